$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 44.533333
$ws.Range("I11").Value = 44.533333
$ws.Range("K11").Value = 44.533333
$ws.Range("M11").Value = 95.466667

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

$ws.Range("H55").Value = 85.61539
$ws.Range("I55").Value = 80.14286
$ws.Range("K55").Value = 80.14286
$ws.Range("M55").Value = 133.85714

$ws.Range("H70").Value = 4058.9
$ws.Range("I70").Value = 1995
$ws.Range("J70").Value = 4288.222
$ws.Range("K70").Value = 5985
$ws.Range("L70").Value = 12864.666
$ws.Range("M70").Value = -5715
$ws.Range("N70").Value = -13404.666

$ws.Range("H73").Value = 4058.9
$ws.Range("I73").Value = 1995
$ws.Range("J73").Value = 4288.222
$ws.Range("K73").Value = 5985
$ws.Range("L73").Value = 12864.666
$ws.Range("M73").Value = -5049
$ws.Range("N73").Value = -14736.666

$ws.Range("H98").Value = 933.4783
$ws.Range("I98").Value = 930.4545000000001
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 930.4545000000001
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 567.5454999999999
$ws.Range("N98").Value = -3996

$ws.Range("H113").Value = 9199.691999999999
$ws.Range("I113").Value = 8319.6
$ws.Range("K113").Value = 8319.6
$ws.Range("M113").Value = -5065.6

$ws.Range("H122").Value = 933.4783
$ws.Range("I122").Value = 930.4545000000001
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2791.3635
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -341.3635000000004
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 250.16667
$ws.Range("I5").Value = 260.2
$ws.Range("K5").Value = 260.2
$ws.Range("M5").Value = -148.2

$ws.Range("H32").Value = 12663545
$ws.Range("I32").Value = 12992313
$ws.Range("K32").Value = 12992313
$ws.Range("M32").Value = -12992026

$ws.Range("H45").Value = 8741
$ws.Range("I45").Value = 6568.3335
$ws.Range("J45").Value = 12000
$ws.Range("K45").Value = 6568.3335
$ws.Range("L45").Value = 12000
$ws.Range("M45").Value = -6191.3335
$ws.Range("N45").Value = -12754

$ws.Range("H61").Value = 1947.1428
$ws.Range("I61").Value = 1947.1428
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1947.1428
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1735.1428
$ws.Range("N61").ClearContents()

$ws.Range("H130").Value = 149996
$ws.Range("J130").Value = 149996
$ws.Range("L130").Value = 149996
$ws.Range("N130").Value = -160036

$ws.Range("H136").Value = 1947.1428
$ws.Range("I136").Value = 1947.1428
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5841.428400000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3291.428400000001
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 250.16667
$ws.Range("I4").Value = 260.2
$ws.Range("K4").Value = 260.2
$ws.Range("M4").Value = -145.2

$ws.Range("H64").Value = 883.55554
$ws.Range("I64").Value = 877.5
$ws.Range("J64").Value = 885.2857
$ws.Range("K64").Value = 877.5
$ws.Range("L64").Value = 885.2857
$ws.Range("M64").Value = -652.5
$ws.Range("N64").Value = -1335.2857

$ws.Range("H67").Value = 883.55554
$ws.Range("I67").Value = 877.5
$ws.Range("J67").Value = 885.2857
$ws.Range("K67").Value = 877.5
$ws.Range("L67").Value = 885.2857
$ws.Range("M67").Value = -97.5
$ws.Range("N67").Value = -2445.2857

$ws.Range("H86").Value = 3035.6667
$ws.Range("I86").Value = 2281.9167
$ws.Range("J86").Value = 4543.1665
$ws.Range("K86").Value = 2281.9167
$ws.Range("L86").Value = 4543.1665
$ws.Range("M86").Value = -1158.9167
$ws.Range("N86").Value = -6789.1665

$ws.Range("H89").Value = 3035.6667
$ws.Range("I89").Value = 2281.9167
$ws.Range("J89").Value = 4543.1665
$ws.Range("K89").Value = 11409.5835
$ws.Range("L89").Value = 22715.8325
$ws.Range("M89").Value = -5793.583500000001
$ws.Range("N89").Value = -33947.8325

$ws.Range("H137").Value = 69971.82000000001
$ws.Range("J137").Value = 69971.82000000001
$ws.Range("L137").Value = 69971.82000000001
$ws.Range("N137").Value = -80171.82000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1820.8085
$ws.Range("I31").Value = 1705.0952
$ws.Range("K31").Value = 1705.0952
$ws.Range("M31").Value = -1410.0952

$ws.Range("H34").Value = 1820.8085
$ws.Range("I34").Value = 1705.0952
$ws.Range("K34").Value = 1705.0952
$ws.Range("M34").Value = -1503.0952

$ws.Range("H69").Value = 8709.4
$ws.Range("I69").Value = 8709.4
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 8709.4
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -7960.4
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 8709.4
$ws.Range("I72").Value = 8709.4
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 26128.2
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -22384.2
$ws.Range("N72").ClearContents()

$ws.Range("H99").Value = 18704920
$ws.Range("J99").Value = 33342668
$ws.Range("L99").Value = 33342668
$ws.Range("N99").Value = -33345664

$ws.Range("H105").Value = 3658.647
$ws.Range("I105").Value = 3742.7273
$ws.Range("K105").Value = 3742.7273
$ws.Range("M105").Value = -1995.7273

$ws.Range("H107").Value = 13483.875
$ws.Range("J107").Value = 26300.75
$ws.Range("L107").Value = 26300.75
$ws.Range("N107").Value = -30140.75

$ws.Range("H126").Value = 18704920
$ws.Range("J126").Value = 33342668
$ws.Range("L126").Value = 100028004
$ws.Range("N126").Value = -100032944

$ws.Range("H132").Value = 1782
$ws.Range("I132").Value = 1855.7084
$ws.Range("J132").Value = 13
$ws.Range("K132").Value = 5567.1252
$ws.Range("L132").Value = 39
$ws.Range("M132").Value = -3037.1252
$ws.Range("N132").Value = -5099

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 91441.82000000001
$ws.Range("I46").Value = 250325
$ws.Range("J46").Value = 651.4286
$ws.Range("K46").Value = 750975
$ws.Range("L46").Value = 1954.2858
$ws.Range("M46").Value = -750884
$ws.Range("N46").Value = -2136.2858

$ws.Range("H97").Value = 398
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()

$ws.Range("H125").Value = 6127
$ws.Range("I125").Value = 5908.75
$ws.Range("K125").Value = 17726.25
$ws.Range("M125").Value = -12806.25

$ws.Range("H129").Value = 1828.8572
$ws.Range("I129").Value = 1065.5714
$ws.Range("J129").Value = 2592.1428
$ws.Range("K129").Value = 3196.7142
$ws.Range("L129").Value = 7776.428400000001
$ws.Range("M129").Value = 1803.2858
$ws.Range("N129").Value = -17776.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 48470
$ws.Range("J42").Value = 48470
$ws.Range("L42").Value = 48470
$ws.Range("N42").Value = -49440

$ws.Range("H115").Value = 48470
$ws.Range("J115").Value = 48470
$ws.Range("L115").Value = 48470
$ws.Range("N115").Value = -50820

$ws.Range("H132").Value = 2991.7144
$ws.Range("I132").Value = 2260.5
$ws.Range("K132").Value = 6781.5
$ws.Range("M132").Value = -4251.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 723.1429000000001
$ws.Range("I35").Value = 212.4
$ws.Range("K35").Value = 212.4
$ws.Range("M35").Value = 123.6

$ws.Range("H68").Value = 3272.0588
$ws.Range("J68").Value = 3498
$ws.Range("L68").Value = 3498
$ws.Range("N68").Value = -4996

$ws.Range("H71").Value = 3272.0588
$ws.Range("J71").Value = 3498
$ws.Range("L71").Value = 17490
$ws.Range("N71").Value = -24978

$ws.Range("H93").Value = 2132.6875
$ws.Range("I93").Value = 2518.0833
$ws.Range("J93").Value = 976.5
$ws.Range("K93").Value = 2518.0833
$ws.Range("L93").Value = 976.5
$ws.Range("M93").Value = -1270.0833
$ws.Range("N93").Value = -3472.5

$ws.Range("H94").Value = 60000
$ws.Range("J94").Value = 60000
$ws.Range("L94").Value = 60000
$ws.Range("N94").Value = -61352

$ws.Range("H100").Value = 49798.24
$ws.Range("I100").Value = 77418.07000000001
$ws.Range("J100").Value = 8368.5
$ws.Range("K100").Value = 77418.07000000001
$ws.Range("L100").Value = 8368.5
$ws.Range("M100").Value = -76877.07000000001
$ws.Range("N100").Value = -9450.5

$ws.Range("H132").Value = 3848.8572
$ws.Range("J132").Value = 8023.3335
$ws.Range("L132").Value = 24070.0005
$ws.Range("N132").Value = -29130.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H132").Value = 2640.2
$ws.Range("I132").Value = 2645.2632
$ws.Range("K132").Value = 7935.7896
$ws.Range("M132").Value = -5405.7896

$ws.Range("H136").Value = 1990
$ws.Range("I136").Value = 1091.6666
$ws.Range("J136").Value = 4300
$ws.Range("K136").Value = 3274.9998
$ws.Range("L136").Value = 12900
$ws.Range("M136").Value = -724.9998000000001
$ws.Range("N136").Value = -18000
